$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Controls": n_sims value 1 -> 100
# -----------------------------------------------------------------
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Range("B2").Value = 100

# -----------------------------------------------------------------
# Sheet "Bins": length bin midpoints shift down by 2 (41->40 ... 99->98)
# and a new row 32 is appended with B32 = 100 (A32 left empty)
# -----------------------------------------------------------------
$wsBins = $wb.Worksheets.Item("Bins")
for ($r = 2; $r -le 31; $r++) {
    $cell = $wsBins.Cells.Item($r, 2)
    $cell.Value = $cell.Value2 - 1
}
$wsBins.Cells.Item(32, 2).Value = 100

# -----------------------------------------------------------------
# Sheet "Selex": add age-based selectivity parameters alongside the
# existing length-based ones (new rows 6-13), and relabel the Notes
# column so it is clear which rows are length-based vs age-based.
# -----------------------------------------------------------------
$wsSelex = $wb.Worksheets.Item("Selex")

# New Par names for the age-based rows (column A), row order first.
$wsSelex.Range("A6").Value = "fish_age_slope_f"
$wsSelex.Range("A7").Value = "fish_age_slope_m"
$wsSelex.Range("A8").Value = "fish_age_midpoint_f"
$wsSelex.Range("A9").Value = "fish_age_midpoint_m"
$wsSelex.Range("A10").Value = "srv_age_slope_f"
$wsSelex.Range("A11").Value = "srv_age_slope_m"
$wsSelex.Range("A12").Value = "srv_age_midpoint_f"
$wsSelex.Range("A13").Value = "srv_age_midpoint_m"

# Age-based Values (column B) for the new rows.
$wsSelex.Range("B6").Value = 0.7
$wsSelex.Range("B7").Value = 0.9
$wsSelex.Range("B8").Value = 4
$wsSelex.Range("B9").Value = 7
$wsSelex.Range("B10").Value = 0.45
$wsSelex.Range("B11").Value = 0.35
$wsSelex.Range("B12").Value = 2
$wsSelex.Range("B13").Value = 5

# Age-based Notes (column C) for the new rows.
$wsSelex.Range("C6").Value = "slope of logistic selectivity (age-based)"
$wsSelex.Range("C7").Value = "slope of logistic selectivity (age-based)"
$wsSelex.Range("C8").Value = "Midpoint of logistic selectivity (age-based)"
$wsSelex.Range("C9").Value = "Midpoint of logistic selectivity (age-based)"
$wsSelex.Range("C10").Value = "slope of logistic selectivity (age-based)"
$wsSelex.Range("C11").Value = "slope of logistic selectivity (age-based)"
$wsSelex.Range("C12").Value = "Midpoint of logistic selectivity (age-based)"
$wsSelex.Range("C13").Value = "Midpoint of logistic selectivity (age-based)"

# Relabel the pre-existing length-based rows to match, and update their
# example values.
$wsSelex.Range("B2").Value = 1
$wsSelex.Range("C2").Value = "slope of logistic selectivity (length-based)"

$wsSelex.Range("B3").Value = 65
$wsSelex.Range("C3").Value = "Midpoint of logistic selectivity (length-based)"

$wsSelex.Range("B4").Value = 0.35
$wsSelex.Range("C4").Value = "slope of logistic selectivity (length-based)"

$wsSelex.Range("B5").Value = 55
$wsSelex.Range("C5").Value = "Midpoint of logistic selectivity (length-based)"

# -----------------------------------------------------------------
# Sheet "Growth_Param": tweak a couple of values / labels
# -----------------------------------------------------------------
$wsGrowth = $wb.Worksheets.Item("Growth_Param")
$wsGrowth.Range("B6").Value = 0.00000945
$wsGrowth.Range("C6").Value = "alpha "
$wsGrowth.Range("A8").Value = 2
$wsGrowth.Range("B8").Value = 2

# -----------------------------------------------------------------
# Window / sheet-view bookkeeping to mirror the author's edit: active
# sheet moves from Growth_Param back to Controls, and a few cell
# selections move around the workbook.
# -----------------------------------------------------------------
$wsMaturity = $wb.Worksheets.Item("Maturity_At_Age")
$wsRecruit = $wb.Worksheets.Item("Recruitment_Mortality")

$wsBins.Range("C31").Select()
$wsGrowth.Range("C8").Select()
$wsSelex.Range("B3").Select()
$wsRecruit.Range("B8").Select()

$wsControls.Range("B3").Select()
$wsControls.Activate()
